# RS_v2.0_distribution.xlsx update
#
# Summary of the change (from the commit):
#   "Fixed bug in saving initialisation file from GUI for a stage-structured
#    population and revised format of GUI initialisation file to be similar
#    to batch initialisation file"
#
# Concretely, on the "Recipients" sheet the row for Pierre Barry / CEFE CNRS
# (previously recorded against version 170419) is removed from its old
# position and re-added at the end of the table against the new version
# 180822, with all of the rows in between shifting up by one. The cell
# comment that was anchored to the "Javier Babi Almenar" row follows that
# row as it shifts. The "Versions" sheet gains a matching new row recording
# version 180822.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Recipients"
$ws2 = $wb.Worksheets.Item(2)   # "Versions"

# --- Preserve & remove the existing cell comment before the rows shift ---
# It is currently anchored at B14 (Javier Babi Almenar's row); once row 11
# is deleted that row becomes row 13, so the comment must move there too.
$comment = $ws1.Comments.Item(1)
$commentText = $comment.Text()
$comment.Delete()

# --- Remove the old Pierre Barry / CEFE CNRS / 170419 row ---
# This shifts rows 12:27 up to 11:26.
$ws1.Rows("11:11").Delete()

# --- Re-insert a blank row just above the trailing blank rows ---
# After the deletion above, the two originally-blank placeholder rows
# (old rows 26 & 27) have shifted up to become rows 25 & 26. Insert a new
# blank row at 25 so that the placeholder rows end up back at 26 & 27,
# matching the original table size, and row 25 is free for the new entry.
$ws1.Rows("25:25").Insert()

# --- Populate the new row 25 with the updated Pierre Barry entry ---
$ws1.Range("A25").Value = "Pierre Barry"
$ws1.Range("B25").Value = "CEFE CNRS"
$ws1.Range("C25").Value = 180822
$ws1.Range("C25").WrapText = $true
$ws1.Range("C25").VerticalAlignment = -4160

# --- Re-anchor the comment to its new row (B14 -> B13) ---
$ws1.Range("B13").AddComment($commentText)

# --- Record the new version on the "Versions" sheet ---
$ws2.Range("A11").Value = 180822
$ws2.Range("A11").WrapText = $true
$ws2.Range("A11").VerticalAlignment = -4160
